$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 18.67841498217124
$ws.Cells.Item(2, 3).Value = 15.58052727591495
$ws.Cells.Item(2, 4).Value = 6.381431055820749
$ws.Cells.Item(2, 5).Value = 11.55663407063
$ws.Cells.Item(2, 6).Value = 44.04471527644677
$ws.Cells.Item(2, 9).Value = 28.41591799105303
$ws.Cells.Item(2, 10).Value = 9.812307771001098
$ws.Cells.Item(2, 13).Value = 19.86841817526783
$ws.Cells.Item(2, 14).Value = 19.74095577190075
$ws.Cells.Item(3, 2).Value = 18.2568108009636
$ws.Cells.Item(3, 3).Value = 15.17686764404504
$ws.Cells.Item(3, 4).Value = 6.380924773044557
$ws.Cells.Item(3, 5).Value = 11.57567309220733
$ws.Cells.Item(3, 6).Value = 43.91941992684593
$ws.Cells.Item(3, 9).Value = 28.43656366824174
$ws.Cells.Item(3, 10).Value = 9.835241316783501
$ws.Cells.Item(3, 13).Value = 19.75657196039054
$ws.Cells.Item(3, 14).Value = 19.8129151788179
$ws.Cells.Item(4, 2).Value = 17.99843605282638
$ws.Cells.Item(4, 3).Value = 14.92853205735597
$ws.Cells.Item(4, 4).Value = 6.380833982313437
$ws.Cells.Item(4, 5).Value = 11.58867990716307
$ws.Cells.Item(4, 6).Value = 43.85473432854575
$ws.Cells.Item(4, 9).Value = 28.45626921228903
$ws.Cells.Item(4, 10).Value = 9.850360458558093
$ws.Cells.Item(4, 13).Value = 19.69225295019283
$ws.Cells.Item(4, 14).Value = 19.85908292750351
$ws.Cells.Item(5, 2).Value = 17.89343537626812
$ws.Cells.Item(5, 3).Value = 14.82737912150115
$ws.Cells.Item(5, 4).Value = 6.380852546114607
$ws.Cells.Item(5, 5).Value = 11.59431169799893
$ws.Cells.Item(5, 6).Value = 43.83146349030731
$ws.Cells.Item(5, 9).Value = 28.46606021128445
$ws.Cells.Item(5, 10).Value = 9.856782919271541
$ws.Cells.Item(5, 13).Value = 19.66715794034369
$ws.Cells.Item(5, 14).Value = 19.87839697488428
$ws.Cells.Item(6, 2).Value = 17.87602246339847
$ws.Cells.Item(6, 3).Value = 14.81059055635252
$ws.Cells.Item(6, 4).Value = 6.38085899099959
$ws.Cells.Item(6, 5).Value = 11.59526687719406
$ws.Cells.Item(6, 6).Value = 43.82778616641914
$ws.Cells.Item(6, 9).Value = 28.46779213812449
$ws.Cells.Item(6, 10).Value = 9.857865154341622
$ws.Cells.Item(6, 13).Value = 19.66305886946208
$ws.Cells.Item(6, 14).Value = 19.88163431458894
$ws.Cells.Item(7, 2).Value = 17.99701858504246
$ws.Cells.Item(7, 3).Value = 14.92716746110034
$ws.Cells.Item(7, 4).Value = 6.380834007388684
$ws.Cells.Item(7, 5).Value = 11.58875451722235
$ws.Cells.Item(7, 6).Value = 43.85440797262257
$ws.Cells.Item(7, 9).Value = 28.45639413718208
$ws.Cells.Item(7, 10).Value = 9.850446015758282
$ws.Cells.Item(7, 13).Value = 19.69190996784488
$ws.Cells.Item(7, 14).Value = 19.85934137605479
$ws.Cells.Item(8, 2).Value = 18.53304174297687
$ws.Cells.Item(8, 3).Value = 15.44154547841061
$ws.Cells.Item(8, 4).Value = 6.381210946943814
$ws.Cells.Item(8, 5).Value = 11.56292560674216
$ws.Cells.Item(8, 6).Value = 43.99897862159669
$ws.Cells.Item(8, 9).Value = 28.42157396885947
$ws.Cells.Item(8, 10).Value = 9.820000038012028
$ws.Cells.Item(8, 13).Value = 19.82896247084922
$ws.Cells.Item(8, 14).Value = 19.76535631902215
$ws.Cells.Item(9, 2).Value = 19.58105183694634
$ws.Cells.Item(9, 3).Value = 16.43933353184698
$ws.Cells.Item(9, 4).Value = 6.383687001542428
$ws.Cells.Item(9, 5).Value = 11.52271041742183
$ws.Cells.Item(9, 6).Value = 44.37906107582296
$ws.Cells.Item(9, 9).Value = 28.40932922777306
$ws.Cells.Item(9, 10).Value = 9.768516665222132
$ws.Cells.Item(9, 13).Value = 20.13130332436374
$ws.Cells.Item(9, 14).Value = 19.59673510657795
$ws.Cells.Item(10, 2).Value = 20.3401327206374
$ws.Cells.Item(10, 3).Value = 17.15691700166241
$ws.Cells.Item(10, 4).Value = 6.386553371959593
$ws.Cells.Item(10, 5).Value = 11.49950875631154
$ws.Cells.Item(10, 6).Value = 44.71616186142779
$ws.Cells.Item(10, 9).Value = 28.43480893340827
$ws.Cells.Item(10, 10).Value = 9.735684931718282
$ws.Cells.Item(10, 13).Value = 20.37249784705762
$ws.Cells.Item(10, 14).Value = 19.48232119029392
$ws.Cells.Item(11, 2).Value = 20.68134180755711
$ws.Cells.Item(11, 3).Value = 17.4783183209747
$ws.Cells.Item(11, 4).Value = 6.388082347605699
$ws.Cells.Item(11, 5).Value = 11.49032751390431
$ws.Cells.Item(11, 6).Value = 44.8817872967395
$ws.Cells.Item(11, 9).Value = 28.45393521285749
$ws.Cells.Item(11, 10).Value = 9.721829280286682
$ws.Cells.Item(11, 13).Value = 20.48603514792891
$ws.Cells.Item(11, 14).Value = 19.43230925413039
$ws.Cells.Item(12, 2).Value = 20.80982493380829
$ws.Cells.Item(12, 3).Value = 17.59917485505603
$ws.Cells.Item(12, 4).Value = 6.388693479177836
$ws.Cells.Item(12, 5).Value = 11.48704795270112
$ws.Cells.Item(12, 6).Value = 44.94624041318179
$ws.Cells.Item(12, 9).Value = 28.4622637386335
$ws.Cells.Item(12, 10).Value = 9.716737480519392
$ws.Cells.Item(12, 13).Value = 20.52954822682931
$ws.Cells.Item(12, 14).Value = 19.41366238944053
$ws.Cells.Item(13, 2).Value = 20.78218811574155
$ws.Cells.Item(13, 3).Value = 17.57318601664479
$ws.Cells.Item(13, 4).Value = 6.388560435046796
$ws.Cells.Item(13, 5).Value = 11.48774550043745
$ws.Cells.Item(13, 6).Value = 44.93228267346124
$ws.Cells.Item(13, 9).Value = 28.46042172434605
$ws.Cells.Item(13, 10).Value = 9.717827198932959
$ws.Cells.Item(13, 13).Value = 20.52015431093047
$ws.Cells.Item(13, 14).Value = 19.41766537558619
$ws.Cells.Item(14, 2).Value = 20.69192747305264
$ws.Cells.Item(14, 3).Value = 17.48827901391666
$ws.Cells.Item(14, 4).Value = 6.388131982856588
$ws.Cells.Item(14, 5).Value = 11.49005375252591
$ws.Cells.Item(14, 6).Value = 44.88705532346396
$ws.Cells.Item(14, 9).Value = 28.45459863496422
$ws.Cells.Item(14, 10).Value = 9.721407268969147
$ws.Cells.Item(14, 13).Value = 20.48960476863278
$ws.Cells.Item(14, 14).Value = 19.43076932829207
$ws.Cells.Item(15, 2).Value = 20.63654189212788
$ws.Cells.Item(15, 3).Value = 17.43615656340642
$ws.Cells.Item(15, 4).Value = 6.387873722804346
$ws.Cells.Item(15, 5).Value = 11.49149329242959
$ws.Cells.Item(15, 6).Value = 44.85957715446197
$ws.Cells.Item(15, 9).Value = 28.45117327976012
$ws.Cells.Item(15, 10).Value = 9.723620348901212
$ws.Cells.Item(15, 13).Value = 20.47095897425901
$ws.Cells.Item(15, 14).Value = 19.43883381533814
$ws.Cells.Item(16, 2).Value = 20.31773958030697
$ws.Cells.Item(16, 3).Value = 17.1358002768197
$ws.Cells.Item(16, 4).Value = 6.386457959528006
$ws.Cells.Item(16, 5).Value = 11.50013638094057
$ws.Cells.Item(16, 6).Value = 44.70558231144084
$ws.Cells.Item(16, 9).Value = 28.43371083143728
$ws.Cells.Item(16, 10).Value = 9.736612137844425
$ws.Cells.Item(16, 13).Value = 20.36515218776518
$ws.Cells.Item(16, 14).Value = 19.48563045036971
$ws.Cells.Item(17, 2).Value = 20.12101496726605
$ws.Cells.Item(17, 3).Value = 16.95015867783361
$ws.Cells.Item(17, 4).Value = 6.385646918140433
$ws.Cells.Item(17, 5).Value = 11.50579016170808
$ws.Cells.Item(17, 6).Value = 44.61423479003621
$ws.Cells.Item(17, 9).Value = 28.42493008877734
$ws.Cells.Item(17, 10).Value = 9.744858549875991
$ws.Cells.Item(17, 13).Value = 20.30119975860078
$ws.Cells.Item(17, 14).Value = 19.51485920285066
$ws.Cells.Item(18, 2).Value = 20.0074857441555
$ws.Cells.Item(18, 3).Value = 16.84291597377007
$ws.Cells.Item(18, 4).Value = 6.385201624387564
$ws.Cells.Item(18, 5).Value = 11.50917134845457
$ws.Cells.Item(18, 6).Value = 44.56285171452845
$ws.Cells.Item(18, 9).Value = 28.42058883696771
$ws.Cells.Item(18, 10).Value = 9.749703304319446
$ws.Cells.Item(18, 13).Value = 20.26477711321579
$ws.Cells.Item(18, 14).Value = 19.53186246851655
$ws.Cells.Item(19, 2).Value = 19.96898598941277
$ws.Cells.Item(19, 3).Value = 16.80652936634126
$ws.Cells.Item(19, 4).Value = 6.385054502706441
$ws.Cells.Item(19, 5).Value = 11.51033837343512
$ws.Cells.Item(19, 6).Value = 44.54565399765293
$ws.Cells.Item(19, 9).Value = 28.41924068517354
$ws.Cells.Item(19, 10).Value = 9.751361117873028
$ws.Cells.Item(19, 13).Value = 20.2525079448514
$ws.Cells.Item(19, 14).Value = 19.53765243486258
$ws.Cells.Item(20, 2).Value = 20.14199681631574
$ws.Cells.Item(20, 3).Value = 16.96996978915967
$ws.Cells.Item(20, 4).Value = 6.385731062794293
$ws.Cells.Item(20, 5).Value = 11.50517492909915
$ws.Cells.Item(20, 6).Value = 44.62383929629181
$ws.Cells.Item(20, 9).Value = 28.42579139514864
$ws.Cells.Item(20, 10).Value = 9.743970187880263
$ws.Cells.Item(20, 13).Value = 20.30797044469935
$ws.Cells.Item(20, 14).Value = 19.51172792444154
$ws.Cells.Item(21, 2).Value = 20.71845993298342
$ws.Cells.Item(21, 3).Value = 17.51324232538442
$ws.Cells.Item(21, 4).Value = 6.388256959084396
$ws.Cells.Item(21, 5).Value = 11.48937041413472
$ws.Cells.Item(21, 6).Value = 44.90029288570891
$ws.Cells.Item(21, 9).Value = 28.45627953571508
$ws.Cells.Item(21, 10).Value = 9.720351509744358
$ws.Cells.Item(21, 13).Value = 20.49856406776009
$ws.Cells.Item(21, 14).Value = 19.42691247555996
$ws.Cells.Item(22, 2).Value = 21.0909276447912
$ws.Cells.Item(22, 3).Value = 17.86328917262339
$ws.Cells.Item(22, 4).Value = 6.3900950234117
$ws.Cells.Item(22, 5).Value = 11.48019039074772
$ws.Cells.Item(22, 6).Value = 45.09106381564566
$ws.Cells.Item(22, 9).Value = 28.48253462276025
$ws.Cells.Item(22, 10).Value = 9.705818901391467
$ws.Cells.Item(22, 13).Value = 20.62613848218912
$ws.Cells.Item(22, 14).Value = 19.37317965758852
$ws.Cells.Item(23, 2).Value = 20.89256905257503
$ws.Cells.Item(23, 3).Value = 17.67696078624677
$ws.Cells.Item(23, 4).Value = 6.389096953355035
$ws.Cells.Item(23, 5).Value = 11.48498490001611
$ws.Cells.Item(23, 6).Value = 44.98833338792314
$ws.Cells.Item(23, 9).Value = 28.46794217101646
$ws.Cells.Item(23, 10).Value = 9.713492623053822
$ws.Cells.Item(23, 13).Value = 20.5577842580167
$ws.Cells.Item(23, 14).Value = 19.40170278706542
$ws.Cells.Item(24, 2).Value = 20.13251225014275
$ws.Cells.Item(24, 3).Value = 16.96101477964791
$ws.Cells.Item(24, 4).Value = 6.385692955616977
$ws.Cells.Item(24, 5).Value = 11.5054526684728
$ws.Cells.Item(24, 6).Value = 44.6194935650699
$ws.Cells.Item(24, 9).Value = 28.42539979683891
$ws.Cells.Item(24, 10).Value = 9.744371493063101
$ws.Cells.Item(24, 13).Value = 20.30490834235048
$ws.Cells.Item(24, 14).Value = 19.51314295465331
$ws.Cells.Item(25, 2).Value = 19.29882785535307
$ws.Cells.Item(25, 3).Value = 16.17152155471333
$ws.Cells.Item(25, 4).Value = 6.382832423901003
$ws.Cells.Item(25, 5).Value = 11.53247424787545
$ws.Cells.Item(25, 6).Value = 44.26599034030254
$ws.Cells.Item(25, 9).Value = 28.40660658895726
$ws.Cells.Item(25, 10).Value = 9.781566112807704
$ws.Cells.Item(25, 13).Value = 20.04605660915322
$ws.Cells.Item(25, 14).Value = 19.64068186629105

Write-Host "Updated loading_percent values for case with 380 kV"
